$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Sema4d", "Plxnb1", "ECs", 3, 1, 68.50681433333334, 205.520443, 0.9663865053086182, 0.9663865053086185, 2, 0.6666666666666666, 0.448495, 1.345485, 0.05033609585745587, 0.05033609585745586, 30.72496369442834, 276.524673249855, 0.04864412376656639, 0.0486441237665664),
    @("ECs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 68.50681433333334, 205.520443, 0.9663865053086182, 0.9663865053086185, 3, 1, 0.6116653333333333, 1.834996, 0.06864924882406573, 0.06864924882406571, 41.90324342480311, 377.129190823228, 0.06634170766315065, 0.06634170766315065),
    @("ECs", "Sema4d", "Plxnb1", "sCs", 3, 1, 68.50681433333334, 205.520443, 0.9663865053086182, 0.9663865053086185, 3, 1, 7.849847333333334, 23.549542, 0.8810146553184784, 0.8810146553184784, 537.7680338096785, 4839.912304287106, 0.8514006738789013, 0.8514006738789014),
    @("FAPs", "Sema4d", "Plxnb1", "ECs", 3, 1, 1.228643, 3.685929, 0.01733176511849816, 0.01733176511849817, 2, 0.6666666666666666, 0.448495, 1.345485, 0.05033609585745587, 0.05033609585745586, 0.551040242285, 4.959362180565, 0.0008724133903836335, 0.0008724133903836335),
    @("FAPs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 1.228643, 3.685929, 0.01733176511849816, 0.01733176511849817, 3, 1, 0.6116653333333333, 1.834996, 0.06864924882406573, 0.06864924882406571, 0.7515183301426666, 6.763664971283999, 0.001189812656180043, 0.001189812656180043),
    @("FAPs", "Sema4d", "Plxnb1", "sCs", 3, 1, 1.228643, 3.685929, 0.01733176511849816, 0.01733176511849817, 3, 1, 7.849847333333334, 23.549542, 0.8810146553184784, 0.8810146553184784, 9.644659977168667, 86.801939794518, 0.01526953907193448, 0.01526953907193449),
    @("sCs", "Sema4d", "Plxnb1", "ECs", 3, 1, 1.154206333333333, 3.462619, 0.01628172957288352, 0.01628172957288353, 2, 0.6666666666666666, 0.448495, 1.345485, 0.05033609585745587, 0.05033609585745586, 0.5176557694683335, 4.658901925215, 0.0008195587005058391, 0.0008195587005058391),
    @("sCs", "Sema4d", "Plxnb1", "FAPs", 3, 1, 1.154206333333333, 3.462619, 0.01628172957288352, 0.01628172957288353, 3, 1, 0.6116653333333333, 1.834996, 0.06864924882406573, 0.06864924882406571, 0.7059880016137778, 6.353892014524, 0.00111772850473503, 0.00111772850473503),
    @("sCs", "Sema4d", "Plxnb1", "sCs", 3, 1, 1.154206333333333, 3.462619, 0.01628172957288352, 0.01628172957288353, 3, 1, 7.849847333333334, 23.549542, 0.8810146553184784, 0.8810146553184784, 9.060343507833112, 81.543091570498, 0.01434444236764265, 0.01434444236764266)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

Write-Host "Updated rows 2-10"
